# Recalculated results after fixing the "average over return periods" logic so
# it behaves correctly for both single-hazard and multi-hazard provinces.
# Only the derived "Socio-economic capacity" (I) and "Risk to well-being" (J)
# columns move; the underlying inputs in columns A-H are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 75.1462663191585
$ws.Range("J2").Value = 1.55722479629847
$ws.Range("I3").Value = 79.6216285207426
$ws.Range("J3").Value = 3.31158158730186
$ws.Range("I4").Value = 56.1052257849643
$ws.Range("J4").Value = 2.08007137877504
$ws.Range("I5").Value = 71.2337330660548
$ws.Range("J5").Value = 0.0486352888232896
$ws.Range("I6").Value = 141.854641333976
$ws.Range("J6").Value = 0.013521079137468
$ws.Range("I7").Value = 63.1062038401126
$ws.Range("J7").Value = 0.291789609270009
$ws.Range("I8").Value = 188.219825891876
$ws.Range("J8").Value = 0.0595532618371669
$ws.Range("I9").Value = 101.551650880561
$ws.Range("J9").Value = 0.822118526820828
$ws.Range("I10").Value = 66.960752521907
$ws.Range("J10").Value = 0.839752196831348
$ws.Range("I11").Value = 96.3562346117991
$ws.Range("J11").Value = 0.802857133594821
$ws.Range("I12").Value = 71.8956773591612
$ws.Range("J12").Value = 0.275107837852916
$ws.Range("I13").Value = 81.4273712488175
$ws.Range("J13").Value = 0.513070314338659
$ws.Range("I14").Value = 137.709128756234
$ws.Range("J14").Value = 0.0874051221419011
$ws.Range("I15").Value = 110.467077747115
$ws.Range("J15").Value = 0.110255294839778
$ws.Range("I16").Value = 92.8438539549242
$ws.Range("J16").Value = 0.705884081587664
$ws.Range("I17").Value = 106.457886299716
$ws.Range("J17").Value = 0.551305386925744
$ws.Range("I18").Value = 131.252765187473
$ws.Range("J18").Value = 0.0423231974618751
$ws.Range("I19").Value = 192.534034632102
$ws.Range("J19").Value = 0.0781640956534163
$ws.Range("I20").Value = 70.1490824051907
$ws.Range("J20").Value = 0.133555728224807
$ws.Range("I21").Value = 39.6234039312315
$ws.Range("J21").Value = 0.92304693261887
$ws.Range("I22").Value = 47.7808316406414
$ws.Range("J22").Value = 1.88220401417163
$ws.Range("I23").Value = 74.6798802362905
$ws.Range("J23").Value = 0.350695129782999
$ws.Range("I24").Value = 55.3966597621795
$ws.Range("J24").Value = 1.78080691366611
$ws.Range("I25").Value = 100.485010215626
$ws.Range("J25").Value = 1.07200217679924
$ws.Range("I26").Value = 123.152754651869
$ws.Range("J26").Value = 0.780365642080336
$ws.Range("I27").Value = 172.250682255626
$ws.Range("J27").Value = 0.529440235090973
$ws.Range("I28").Value = 110.95202401016
$ws.Range("J28").Value = 0.554493387539422
$ws.Range("I29").Value = 106.190546551591
$ws.Range("J29").Value = 0.0452342512962872
$ws.Range("I30").Value = 239.660419263113
$ws.Range("J30").Value = 0.0436394738608865
$ws.Range("I31").Value = 59.029802043471
$ws.Range("J31").Value = 0.936146408011064
$ws.Range("I32").Value = 46.8949546215255
$ws.Range("J32").Value = 0.00397071981083791
$ws.Range("I33").Value = 93.1114324321849
$ws.Range("J33").Value = 0.1465375673226
$ws.Range("I34").Value = 49.6978416062788
$ws.Range("J34").Value = 1.93261969641378
$ws.Range("I35").Value = 122.191832817456
$ws.Range("J35").Value = 0.640084430727178
$ws.Range("I36").Value = 78.0860652064569
$ws.Range("J36").Value = 0.385579871088775

Write-Output "Updated socio-economic capacity (I) and risk to well-being (J) values for rows 2-36"